$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$row2 = @(1.02, 1.048611145005293, 1.056815874688694, 1.056015028638772, 1.066967185681081, 1.047006889342252, 1.053653978139931, 1.059552549529211, 1.058753901153756, 1.069676312792945, 1.021679010047283)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $row2[$i]
}

$row3 = @(1.02, 1.049613194974895, 1.057469218873087, 1.056866210320243, 1.06780560143377, 1.047241365004617, 1.054304717182346, 1.060020015334928, 1.05941854236296, 1.070330366393827, 1.021899648094838)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3[$i]
}

$row4 = @(1.02, 1.05026214237279, 1.057892396448508, 1.057417793531764, 1.068348892134907, 1.047392189086405, 1.054725725404637, 1.060322229707767, 1.059848774745822, 1.070753711405034, 1.022042290009998)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4[$i]
}

$row5 = @(1.02, 1.050535091676219, 1.058070399391584, 1.057649872117223, 1.068577476704975, 1.047455380138377, 1.054902701152968, 1.060449215600379, 1.060029682910425, 1.070931715244745, 1.022102226187456)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5[$i]
}

$row6 = @(1.02, 1.050580928771253, 1.058100292635854, 1.057688850412077, 1.068615867874651, 1.047465977547003, 1.054932415188621, 1.060470533246718, 1.060060060412789, 1.070961604573721, 1.02211228793614)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $row6[$i]
}

$row7 = @(1.02, 1.050265789017809, 1.057894774545084, 1.057420893820478, 1.068351945767439, 1.04739303429538, 1.054728090227688, 1.060323926756314, 1.05985119189972, 1.070756089786247, 1.022043091000417)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $row7[$i]
}

$row8 = @(1.02, 1.048949677263292, 1.057036587504341, 1.056302520595485, 1.067250369883697, 1.047086316952607, 1.053873911031655, 1.059710586495946, 1.058978484950809, 1.069897325927745, 1.021753601484747)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $row8[$i]
}

$row9 = @(1.02, 1.046634800174194, 1.055527637027563, 1.054338079856863, 1.06531528938022, 1.046538994153726, 1.052368292731427, 1.058627801916892, 1.057441974429183, 1.068385111707588, 1.021242539321749)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $row9[$i]
}

$row10 = @(1.02, 1.045094474165819, 1.054523967933517, 1.053032750718939, 1.064029382786807, 1.046169541716841, 1.051364301595046, 1.057904662394428, 1.056418577076008, 1.067377736686124, 1.020901219721667)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $row10[$i]
}

$row11 = @(1.02, 1.044428199446057, 1.05408993194974, 1.05246856457131, 1.063473573094704, 1.046008487592323, 1.05092951519133, 1.057591243620051, 1.055975671205119, 1.066941728939262, 1.020753284429059)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $row11[$i]
}

$row12 = @(1.02, 1.044180820716586, 1.053928797205173, 1.052259156754139, 1.063267271787045, 1.045948503377842, 1.050768009263723, 1.057474782699449, 1.055811192093957, 1.066779806294438, 1.020698313774813)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $row12[$i]
}

$row13 = @(1.02, 1.044233879514241, 1.053963357259573, 1.052304068359901, 1.063311517253971, 1.04596137749699, 1.050802653131828, 1.057499765902158, 1.055846471783232, 1.066814537873201, 1.020710106106372)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $row13[$i]
}

$row14 = @(1.02, 1.044407748889514, 1.054076610747676, 1.052451251665186, 1.063456517072152, 1.046003532570752, 1.050916165196734, 1.057581617800027, 1.055962074573346, 1.066928343743312, 1.020748740961768)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $row14[$i]
}

$row15 = @(1.02, 1.044514889533592, 1.054146401340185, 1.052541956872605, 1.063545876319776, 1.046029484313132, 1.050986102837858, 1.057632043757679, 1.056033306048199, 1.066998467311915, 1.020772542413304)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Value = $row15[$i]
}

$row16 = @(1.02, 1.045138707361797, 1.054552785386095, 1.053070215706538, 1.06406629116289, 1.046180207665425, 1.05139315592226, 1.05792545684717, 1.056447976252263, 1.067406677234286, 1.020911034747956)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $row16[$i]
}

$row17 = @(1.02, 1.045530198804942, 1.054807850182745, 1.053401855186329, 1.064393001432812, 1.046274463997403, 1.051648476352433, 1.058109428980599, 1.056708150719949, 1.067662788830474, 1.020997869675568)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
}

$row18 = @(1.02, 1.04575861630466, 1.054956679032111, 1.053595394175506, 1.064583662040404, 1.046329337998527, 1.051797395401861, 1.058216708321077, 1.056859928396236, 1.067812192866552, 1.021048505326578)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "18").Value = $row18[$i]
}

$row19 = @(1.02, 1.045836512139976, 1.055007434914834, 1.053661402771494, 1.064648688624923, 1.046348030937654, 1.051848172068791, 1.058253282939822, 1.056911684415396, 1.067863138880419, 1.021065768440816)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
}

$row20 = @(1.02, 1.045488188525824, 1.054780478561171, 1.053366263129193, 1.064357938568867, 1.046264361951053, 1.05162108340148, 1.058089693451462, 1.056680234146823, 1.067635308560685, 1.020988554523211)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
}

$row21 = @(1.02, 1.044356545777685, 1.05404325804463, 1.052407905532395, 1.063413814055182, 1.045991123407266, 1.050882738906371, 1.057557515641407, 1.055928031435284, 1.066894829912399, 1.020737364533668)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "21").Value = $row21[$i]
}

$row22 = @(1.02, 1.043645647046628, 1.053580233929145, 1.051806251506481, 1.062821080790485, 1.045818393125104, 1.050418472931539, 1.057222664706337, 1.055455299649068, 1.06642943600475, 1.020579310687384)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $row22[$i]
}

$row23 = @(1.02, 1.0440224498567, 1.053825644285386, 1.052125113621586, 1.063135216451008, 1.045910049104241, 1.050664592565096, 1.057400198841453, 1.055705883669883, 1.066676133202304, 1.020663109337945)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "23").Value = $row23[$i]
}

$row24 = @(1.02, 1.045507170955131, 1.054792846452494, 1.053382345341105, 1.064373781671248, 1.046268926952731, 1.051633461112529, 1.05809861117572, 1.056692848375715, 1.067647725654759, 1.020992763681939)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "24").Value = $row24[$i]
}

$row25 = @(1.02, 1.047232739290678, 1.055917338415653, 1.054845183074073, 1.065814829579909, 1.046681298264929, 1.052757578394785, 1.058907958386913, 1.057839037550614, 1.068775925152006, 1.021374770652022)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "25").Value = $row25[$i]
}
